$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "location"
$ws.Range("A2").Value = "narender"
$ws.Range("B2").Value = "delhi"

$ws.Range("A3").Select() | Out-Null
